$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '62.051.58'
$ws.Cells.Item(2, 5).Value = '  -0.87%  '
$ws.Cells.Item(3, 4).Value = '3.411.52'
$ws.Cells.Item(3, 5).Value = '  -1.78%  '
$ws.Cells.Item(5, 4).Value = '407.30'
$ws.Cells.Item(5, 5).Value = '  -1.27%  '
$ws.Cells.Item(6, 4).Value = '133.81'
$ws.Cells.Item(6, 5).Value = '  +3.78%  '
$ws.Cells.Item(7, 4).Value = '0.595'
$ws.Cells.Item(7, 5).Value = '  -0.61%  '
$ws.Cells.Item(8, 5).Value = '  -0.07%  '
$ws.Cells.Item(9, 5).Value = '  -1.41%  '
$ws.Cells.Item(10, 5).Value = '  -5.82%  '
$ws.Cells.Item(11, 4).Value = '42.75'
$ws.Cells.Item(11, 5).Value = '  -2.27%  '
$ws.Cells.Item(12, 5).Value = '  -1.16%  '
$ws.Cells.Item(13, 4).Value = '8.43'
$ws.Cells.Item(13, 5).Value = '  -3.81%  '
$ws.Cells.Item(14, 4).Value = '19.87'
$ws.Cells.Item(14, 5).Value = '  -1.88%  '
$ws.Cells.Item(15, 4).Value = '3.435.29'
$ws.Cells.Item(15, 5).Value = '  -0.81%  '
$ws.Cells.Item(16, 4).Value = '62.120.06'
$ws.Cells.Item(16, 5).Value = '  -0.61%  '
$ws.Cells.Item(17, 5).Value = '  -3.22%  '
$ws.Cells.Item(18, 4).Value = '11.01'
$ws.Cells.Item(18, 5).Value = '  -1.07%  '
$ws.Cells.Item(19, 5).Value = '  -4.84%  '
$ws.Cells.Item(20, 5).Value = '  -5.39%  '
$ws.Cells.Item(21, 4).Value = '84.14'
$ws.Cells.Item(21, 5).Value = '  +2.19%  '
$ws.Cells.Item(22, 4).Value = '313.42'
$ws.Cells.Item(22, 5).Value = '  +0.30%  '
$ws.Cells.Item(23, 5).Value = '  -2.71%  '
$ws.Cells.Item(24, 4).Value = '3.16'
$ws.Cells.Item(24, 5).Value = '  -1.02%  '
$ws.Cells.Item(25, 5).Value = '  +9.86%  '
$ws.Cells.Item(26, 4).Value = '29.59'
$ws.Cells.Item(26, 5).Value = '  -2.72%  '
$ws.Cells.Item(27, 4).Value = '8.18'
$ws.Cells.Item(27, 5).Value = '  +0.37%  '
$ws.Cells.Item(28, 5).Value = '  +4.66%  '
$ws.Cells.Item(29, 4).Value = '7.59'
$ws.Cells.Item(29, 5).Value = '  -3.33%  '
$ws.Cells.Item(30, 5).Value = '  -2.34%  '
$ws.Cells.Item(31, 5).Value = '  -4.09%  '
$ws.Cells.Item(32, 4).Value = '42.77'
$ws.Cells.Item(32, 5).Value = '  -4.59%  '
$ws.Cells.Item(33, 5).Value = '  -0.22%  '
$ws.Cells.Item(34, 5).Value = '  -6.36%  '
$ws.Cells.Item(35, 5).Value = '  -2.37%  '
$ws.Cells.Item(36, 4).Value = '51.83'
$ws.Cells.Item(36, 5).Value = '  -1.66%  '
$ws.Cells.Item(37, 4).Value = '0.999'
$ws.Cells.Item(37, 5).Value = '  +0.29%  '
$ws.Cells.Item(38, 5).Value = '  -4.37%  '
$ws.Cells.Item(39, 4).Value = '2.94'
$ws.Cells.Item(39, 5).Value = '  -3.26%  '
$ws.Cells.Item(40, 5).Value = '  -0.21%  '
$ws.Cells.Item(41, 5).Value = '  -0.56%  '
$ws.Cells.Item(42, 4).Value = '137.27'
$ws.Cells.Item(42, 5).Value = '  -0.43%  '
$ws.Cells.Item(43, 5).Value = '  +2.80%  '
$ws.Cells.Item(44, 5).Value = '  +0.55%  '
$ws.Cells.Item(45, 4).Value = '16.75'
$ws.Cells.Item(45, 5).Value = '  -6.55%  '
$ws.Cells.Item(46, 5).Value = '  -2.56%  '
$ws.Cells.Item(47, 4).Value = '21.24'
$ws.Cells.Item(47, 5).Value = '  -5.81%  '
$ws.Cells.Item(48, 4).Value = '2.121.41'
$ws.Cells.Item(48, 5).Value = '  -4.60%  '
$ws.Cells.Item(49, 4).Value = '2.32'
$ws.Cells.Item(50, 5).Value = '  +3.23%  '
$ws.Cells.Item(51, 4).Value = '1.65'
$ws.Cells.Item(51, 5).Value = '  +16.02%  '
